# Trade #42 closed at 2026-02-17 13:27:31 - unknown UNKNOWN +0.000%
#
# Updates the live trading results workbook with the newly-closed trade:
#  - Summary sheet totals (capital, P&L, trade/win counts, win rate)
#  - Strategy Status row for MarketMaking
#  - A new trade row (#42 / sheet row 43) appended to both the
#    "All Trades" and "MarketMaking" trade logs

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1197.61   # Current Capital
$summary.Range("B4").Value = -2.39     # Total P&L $
$summary.Range("B5").Value = -1.14     # Total P&L %
$summary.Range("B6").Value = 42        # Total Trades
$summary.Range("B7").Value = 17        # Winning Trades
$summary.Range("B9").Value = 40.48     # Win Rate %

# ---------------------------------------------------------------------
# 2. Strategy Status sheet (MarketMaking row)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 97.61      # Capital
$status.Range("D4").Value = 42         # Trades
$status.Range("E4").Value = -2.39      # P&L $
$status.Range("F4").Value = -2.39      # P&L %
$status.Range("G4").Value = 40.48      # Win Rate %

# ---------------------------------------------------------------------
# 3. Append the new trade (#42) as sheet row 43 on both trade-log sheets.
#    Copy the previous row so text columns (date/time/strategy/etc.) keep
#    their string type, then overwrite just the cells that changed.
# ---------------------------------------------------------------------
function Add-Trade43($ws) {
    $ws.Range("A42:Q42").Copy()
    $ws.Range("A43:Q43").PasteSpecial()

    $ws.Cells.Item(43, 1).Value = 42            # Trade #
    $ws.Cells.Item(43, 3).Value = "13:27:24"    # Time
    $ws.Cells.Item(43, 6).Value = 0.68          # Entry Price
    $ws.Cells.Item(43, 7).Value = 0.7           # Exit Price
    $ws.Cells.Item(43, 9).Value = 2.9412        # P&L %
    $ws.Cells.Item(43, 10).Value = 0.02         # P&L $
    $ws.Cells.Item(43, 11).Value = 97.61        # Capital After
    # Columns B, D, E, H, L, M, N, O, P, Q are unchanged from row 42.
}

$allTrades = $wb.Worksheets.Item("All Trades")
Add-Trade43 $allTrades

$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-Trade43 $marketMaking
